$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# Change 1
Replace-Text 'Ativação: 01/01/2011' 'Ativação: 01/01/2023'

# Change 2
Replace-Text 'Curso (semestre ideal): EM (4)' 'Curso (semestre ideal): EF (4), EM (4)'

# Change 3
Replace-Text 'O objetivo desta disciplina é apresentar ao aluno fundamentos de Físico-Química, mais especificamente Termodinâmica Química, com foco em problemas encontrados na Engenharia de Materiais. Bastante ênfase é dada nos assuntos relativos a Termoquímica e Equilíbrio Química, abordando equilíbrio com fases condensadas.' 'Esta disciplina visa apresentar fundamentos de termodinâmica aplicada à área de ciência e engenharia de materiais. Especial ênfase é dada à energia na forma e calor para aquecimento de sistemas termodinâmicos; cálculos de variação de entalpia; entropia e energia de Gibbs de elementos e compostos em mudanças de estado; cálculos de variação de entalpia; entropia e energia de Gibbs de reação; aplicação da propriedade  energia de Gibbs para avaliação de transformações espontâneas e em equilíbrio; fundamentos de termodinâmica de soluções; cálculos de condições de equilíbrio em sistemas heterogêneos. Apresenta-se também as principais diferenças entre esta disciplina e a disciplina de Termodinâmica de Máquinas.'

# Change 4
Replace-Text '1) 1a Lei da Termodinâmica2) 2a e 3a Leis da Termodinâmica3) Relações entre Propriedades Termodinâmicas4) Equilíbrio5) Equilíbrio Químico6) Soluções' '1) Introdução; 2) 1a Lei da Termodinâmica 3) 2a e 3a Leis da Termodinâmica 4) Equilíbrio heterogêneo: composição variável da fase gasosa; 5) Equilíbrio heterogêneo: composição variável da fase condensada;'

# Change 5
Replace-Text '1- 1a Lei da Termodinâmica: sistema e vizinhança; transferência de energia; energia de um sistema; energia como uma função de estado; trabalho; sistema fechado; propriedades extensivas e intensivas; sistema aberto; entalpia; estado estacionário; capacidade térmica a volume constante; capacidade térmica a volume constante; equação de estado, gases não-ideais; expansão e compressão adiabática; entalpias de formação; variação de entalpia em reações químicas; variação de temperatura associadas à reações químicas em sistemas adiabáticos.2- 2a e 3a Leis da Termodinâmica: Entropia como função de estado; variações de entropia associadas à variações de temperatura e pressão; variações de entropia em reações químicas; terceiro princípio da termodinâmica.3- Relações entre Propriedades Termodinâmicas: As funções A e G; potencial químico; grandezas parciais molares; relações entre propriedades derivadas de U, H, A e G; gás ideal; entropia de mistura; capacidade térmica; variação de capacidade térmica; Relação P-T isoentrópica; compressão isoentrópica de sólidos.4- Equilíbrio: Condições de equilíbrio; equilíbrio de fases; variação de pressão de equilíbrio com a temperatura; equação de Clapeyron; variação da pressão de vapor de uma fase condensada com a pressão total aplicada; variação da pressão de vapor com tamanho de partícula.5- Equilíbrio Químico: atividade; equilíbrio químico; equilíbrio em fase gasosa; equilíbrio sólido-vapor; fontes de informação em DGo; diagrama de Ellingham; variação da constante de equilíbrio com a temperatura; gases dissolvidos em metais (Lei de Sievert); equilíbrio químico e temperatura adiabática de chama.6- Soluções: grandezas parciais molares relativas; entropia de mistura - solução ideal; entalpia de mistura  solução ideal; Soluções não-ideais; relação de Gibbs-Duhem; soluções regulares.' '1- Introdução: sistema; vizinhanças; fases; equilíbrio; fronteiras adiabáticas e diatérmicas; processos reversíveis e irreversíveis; estado termodinâmico; mudança de estado; processos cíclicos; equação de estado; calor; trabalho.2- A 1ª lei de Termodinâmica: energia interna; capacidades térmicas; entalpia; entalpia de transformação de fases; entalpia de formação e de reação; entalpia de reação em função da temperatura (introdução ao loop termodinâmico).3- A 2ª e 3ª leis da Termodinâmica: Dispersão de energia e entropia; entropia no zero absoluto; entropia de reação; entropia de reação em função da temperatura; desigualdade de Clausius; critérios de espontaneidade e equilíbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de reação em função da temperatura; equação de Gibbs-Helmholtz.4- Equilíbrio heterogêneo: composição variável da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um gás ideal; pressão de equilíbrio em sistemas metal-óxido-O2(g).5- Equilíbrio heterogêneo: composição variável da fase condensada: fugacidade; atividade termodinâmica; soluções e grandezas parciais molares; potencial químico; modelos de soluções; propriedades termodinâmicas de excesso'

# Change 6
Replace-Text 'O curso será ministrado na forma de aulas expositivas.' 'Esta é uma disciplina fundamental, exigindo dedicação individual para assimilação de definições e conceitos. Isto envolve leitura concentrada e realização de exercícios numéricos.'

# Change 7
Replace-Text 'Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R), com pontuação de 0 a 10, que levará ao cálculo da média final (MF) através da seguinte expressão:MF=(NF+R)/2' 'Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# Change 8
Replace-Text 'REFERÊNCIAS BIBLIOGRÁFICAS1) Castellan, G., Fundamentos de Físico-Química, vol.1, Livros Técnicos e Científicos Editora S.A., 1986.2) Atkins, P.W., Physical Chemistry, Oxford University Press, 1994.3) Pilla, L., Físico Química, vol.1, Livros Técnicos e Científicos Editora S.A., 1979.4) Moore, W.J., Físico Química, vol.1, Editora Edgard Blücher Ltda, 1976.5) Darken, L. & Gurry, R., Physical Chemistry f Metals, McGraw-Hill Book Company Inc., 1953.6) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995.7) Ragone, D.V., Thermodynamics of Materials, vol.1, John Wiley & Sons Inc., 1995. (Livro Texto)8) Swalin, R.A. Thermodynamics of Solids, John Wiley & Sons, 1972.9) Shoemaker, D.P., Garland, C.W., Nibler, J.W., Experiments in Physical Chemistry, McGraw-Hill Book Company, 1989.' '1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. Físico-Química, Livros Técnicos e Científicos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.'

# Change 9
Replace-Text 'LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)' 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)'
